$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.82"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.09"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.345"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05942"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.393"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8183"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9611"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1428"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03533"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07396"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03036"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09407"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.000"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001585"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04803"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "One"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005914"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006265"

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004142"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009867"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "20BitKanKAN"

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00009706"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.742"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.164"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "23BTSETokenBTSE"

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3268"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1333"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "25ProBitTokenPROB"

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002463"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "26UpBotsUBXT"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03926"

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1075"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002702"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003044"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005350"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005299"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04338"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "47BOLOBOLO"
